# Experiment order generation script re-run:
# regenerates randomized task-order CSV filenames / values for each
# task-order sheet, and reorders the RS sheet to be right after GNG.

$wb = $excel.ActiveWorkbook

# --- Step 1: Reorder sheets -------------------------------------------
# Before: GNG, NB, RS, TOL, vSAT
# After:  GNG, RS, NB, TOL, vSAT
# Move RS (currently 3rd) to just before NB (currently 2nd).
$wsRS = $wb.Worksheets.Item("RS_TO-16512556253103416")
$wsNB = $wb.Worksheets.Item("NB_TO-16512556253085067")
$wsRS.Move($wsNB)

# --- Step 2: Rename sheets to the newly generated task-order ids -------
$wb.Worksheets.Item("GNG_TO-16512556209355602").Name = "GNG_TO-16515890418614779"
$wb.Worksheets.Item("RS_TO-16512556253103416").Name  = "RS_TO-16515890418614779"
$wb.Worksheets.Item("NB_TO-16512556253085067").Name  = "NB_TO-16515890429921908"
$wb.Worksheets.Item("TOL_TO-1651255625356043").Name  = "TOL_TO-16515890430390368"
$wb.Worksheets.Item("vSAT_TO-1651255625433689").Name = "vSAT_TO-16515890431015334"

# --- Step 3: Update GNG sheet stim file order ---------------------------
$wsGNG = $wb.Worksheets.Item("GNG_TO-16515890418614779")
$wsGNG.Range("B2").Value = "go_stims-1651589041830263.csv"
$wsGNG.Range("B3").Value = "GNG_stims-16515890418458896.csv"
$wsGNG.Range("B4").Value = "go_stims-16515890418458896.csv"
$wsGNG.Range("B5").Value = "GNG_stims-16515890418614779.csv"

# --- Step 4: Update RS sheet (eyes open / eyes closed order) -----------
$wsRS2 = $wb.Worksheets.Item("RS_TO-16515890418614779")
$wsRS2.Range("B2").Value = "eyes closed"
$wsRS2.Range("B3").Value = "eyes open"

# --- Step 5: Update NB sheet stim file order ----------------------------
$wsNB2 = $wb.Worksheets.Item("NB_TO-16515890429921908")
$wsNB2.Range("B2").Value  = "TB-16515890429765658.csv"
$wsNB2.Range("B3").Value  = "ZB-match_7-1651589042086581.csv"
$wsNB2.Range("B4").Value  = "OB-16515890422780278.csv"
$wsNB2.Range("B5").Value  = "TB-16515890428014767.csv"
$wsNB2.Range("B6").Value  = "ZB-match_2-16515890421803281.csv"
$wsNB2.Range("B7").Value  = "OB-165158904225842.csv"
$wsNB2.Range("B8").Value  = "ZB-match_6-1651589042242828.csv"
$wsNB2.Range("B9").Value  = "OB-16515890426743877.csv"
$wsNB2.Range("B10").Value = "TB-1651589042754635.csv"

# --- Step 6: Update TOL sheet stim file order ---------------------------
$wsTOL = $wb.Worksheets.Item("TOL_TO-16515890430390368")
$wsTOL.Range("B2").Value = "MM_stims-1651589043007842.csv"
$wsTOL.Range("B3").Value = "ZM_stims-16515890429921908.csv"
$wsTOL.Range("B4").Value = "MM_stims-16515890430234113.csv"
$wsTOL.Range("B5").Value = "ZM_stims-1651589043007842.csv"
$wsTOL.Range("B6").Value = "MM_stims-16515890430390368.csv"
$wsTOL.Range("B7").Value = "ZM_stims-16515890430234113.csv"

# --- Step 7: Update vSAT sheet stim file order ---------------------------
$wsVSAT = $wb.Worksheets.Item("vSAT_TO-16515890431015334")
$wsVSAT.Range("B2").Value = "vSAT_stims-16515890430702853.csv"
$wsVSAT.Range("B3").Value = "SAT_stims-16515890430546608.csv"
$wsVSAT.Range("B4").Value = "SAT_stims-16515890430390368.csv"
$wsVSAT.Range("B5").Value = "vSAT_stims-1651589043085909.csv"
